# Update project dates (Sprint Backlog 2 + burn-down chart data) and
# switch the date columns to a date-only (no time) display format.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sprint Backlog 2": column G holds "ESTIMATED DATE OF TASK COMPLETION"
# for data rows 4..93. Re-date them (new run starts 2025-09-08, serial 45908)
# and switch their display format from "yyyy-mm-dd h:mm:ss" to "yyyy-mm-dd".
# ---------------------------------------------------------------------------
$wsBacklog = $wb.Worksheets.Item("Sprint Backlog 2")

$wsBacklog.Range("G4:G93").NumberFormat = "yyyy-mm-dd"

for ($row = 4; $row -le 93; $row++) {
    $newDate = 45908 + [Math]::Floor(($row - 4) / 3)
    $wsBacklog.Cells.Item($row, 7).Value = $newDate
}

# ---------------------------------------------------------------------------
# Sheet "SB BD Ch02": burn-down chart data. Column A = date, B = planned
# total tasks remaining, C = actual incomplete tasks. The sprint grew from
# 31 days (rows 3..33) to 38 days (rows 3..40), starting 2025-09-08.
# ---------------------------------------------------------------------------
$wsChart = $wb.Worksheets.Item("SB BD Ch02")

$wsChart.Range("A3:A40").NumberFormat = "yyyy-mm-dd"

for ($row = 3; $row -le 40; $row++) {
    $newDate = 45908 + ($row - 3)
    $wsChart.Cells.Item($row, 1).Value = $newDate

    $remaining = [Math]::Round(90 * (40 - $row) / 37, 1)
    $wsChart.Cells.Item($row, 2).Value = $remaining
    $wsChart.Cells.Item($row, 3).Value = $remaining
}
